$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the header formatting (bold, centered, bordered - style index 1)
# from an existing header cell onto the new A1 header cell, then set its
# text to "Category".
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "Category"

# The category label cells A2:A46 previously used the same header style
# as row 1, but should now match the plain (unstyled) formatting of the
# rest of the data cells (e.g. B2), so copy that formatting over them.
$ws.Range("B2").Copy()
$ws.Range("A2:A46").PasteSpecial(-4122)
